# Generate Report for Handback
# Adds the handback row for file 61e6f3ee-da59-4ee6-9f63-994c5cc74773 to the
# Overview / zh-cn / de-de sheets (new row 4 on each), mirroring the existing
# rows for cf09f27f-5f85-487d-86be-cbbc8cf4cb3a (which is "in sync with en-US").

$wb = $excel.ActiveWorkbook

# Cornflower-blue underlined color used by the workbook's "HyperLink" cell
# style (RGB 0x64,0x95,0xED == hex FF6495ED), expressed as an Excel BGR color
# value (R + G*256 + B*65536) for Font.Color.
$hyperlinkColor = 15570276

function Set-HandbackHyperlink($ws, $cellRef, $url, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText) | Out-Null
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

$guid = "61e6f3ee-da59-4ee6-9f63-994c5cc74773"
$mdName = "$guid.md"
$zhXlfName = "$guid.2eabd2d5ee7f217fc79ce8dff186292f0a3f2132.zh-cn.xlf"
$deXlfName = "$guid.2eabd2d5ee7f217fc79ce8dff186292f0a3f2132.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

Set-HandbackHyperlink $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/e2e/$mdName" $mdName

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Range("C4").Value = $zhXlfName
$wsZh.Range("D4").Value = "2016-02-22 17:25:21"
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $zhXlfName
$wsZh.Range("G4").Value = "2016-02-22 17:26:12"
$wsZh.Range("H4").Value = "Include"

$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Set-HandbackHyperlink $wsZh "A4" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/e2e/$mdName" $mdName
Set-HandbackHyperlink $wsZh "C4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName" $zhXlfName
Set-HandbackHyperlink $wsZh "E4" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/e2e/$mdName" $mdName
Set-HandbackHyperlink $wsZh "F4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName" $zhXlfName

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Range("C4").Value = $deXlfName
$wsDe.Range("D4").Value = "2016-02-22 17:25:32"
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $deXlfName
$wsDe.Range("G4").Value = "2016-02-22 17:26:32"
$wsDe.Range("H4").Value = "Include"

$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Set-HandbackHyperlink $wsDe "A4" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/e2e/$mdName" $mdName
Set-HandbackHyperlink $wsDe "C4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName" $deXlfName
Set-HandbackHyperlink $wsDe "E4" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/e2e/$mdName" $mdName
Set-HandbackHyperlink $wsDe "F4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61e6f3ee-da59-4ee6-9f63-994c5cc74773/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName" $deXlfName
